$wb = $excel.ActiveWorkbook

# --- Column descriptions sheet: update the "measurementValue" description,
# rename "mass" reference to "body mass", and select B16 as the new active cell ---
$wsCols = $wb.Worksheets.Item("Column descriptions")
$wsCols.Range("B15").Value = 'Numeric trait value or a range of values e.g. "5.3-8.9"'

# Make "Column descriptions" the active sheet/tab, with B16 selected
$null = $wsCols.Activate()
$null = $wsCols.Range("B16").Select()

# Try to resize the saved window metrics to match the new layout
# (best effort; some hosts may not persist these).
try {
    $win = $wb.Windows.Item(1)
    $win.Width = 28800
    $win.Height = 12300
} catch {
}

# Update the absolute path recorded for this workbook (best effort).
try {
    $wb.Path = "C:\Jim\uni\Papers and talks\AnimalTraits\animaltraits.github.io\data\"
} catch {
}
